$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4022.2222
$ws.Range("I76").Value = 4250
$ws.Range("J76").Value = 3840
$ws.Range("K76").Value = 4250
$ws.Range("L76").Value = 3840
$ws.Range("M76").Value = -3935
$ws.Range("N76").Value = -4470
$ws.Range("H79").Value = 4022.2222
$ws.Range("I79").Value = 4250
$ws.Range("J79").Value = 3840
$ws.Range("K79").Value = 4250
$ws.Range("L79").Value = 3840
$ws.Range("M79").Value = -3158
$ws.Range("N79").Value = -6024
$ws.Range("H86").Value = 2108.8235
$ws.Range("I86").Value = 2153.8333
$ws.Range("J86").Value = 2000.8
$ws.Range("K86").Value = 2153.8333
$ws.Range("L86").Value = 2000.8
$ws.Range("M86").Value = -1030.8333
$ws.Range("N86").Value = -4246.8
$ws.Range("H89").Value = 2108.8235
$ws.Range("I89").Value = 2153.8333
$ws.Range("J89").Value = 2000.8
$ws.Range("K89").Value = 10769.1665
$ws.Range("L89").Value = 10004
$ws.Range("M89").Value = -5153.166499999999
$ws.Range("N89").Value = -21236
$ws.Range("H137").Value = 1508.15
$ws.Range("I137").Value = 1506.2609
$ws.Range("J137").Value = 1510.7059
$ws.Range("K137").Value = 4518.7827
$ws.Range("L137").Value = 4532.1177
$ws.Range("M137").Value = -1968.7827
$ws.Range("N137").Value = -9632.117699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 846.8200000000001
$ws.Range("I32").Value = 779.1123700000001
$ws.Range("J32").Value = 1394.6364
$ws.Range("K32").Value = 779.1123700000001
$ws.Range("L32").Value = 1394.6364
$ws.Range("M32").Value = -492.1123700000001
$ws.Range("N32").Value = -1968.6364
$ws.Range("H74").Value = 810.1389
$ws.Range("I74").Value = 750.4545000000001
$ws.Range("K74").Value = 750.4545000000001
$ws.Range("M74").Value = 123.5454999999999
$ws.Range("H77").Value = 810.1389
$ws.Range("I77").Value = 750.4545000000001
$ws.Range("K77").Value = 3752.2725
$ws.Range("M77").Value = 615.7275
$ws.Range("H132").Value = 1432.4717
$ws.Range("I132").Value = 960.6286
$ws.Range("J132").Value = 2349.9443
$ws.Range("K132").Value = 2881.8858
$ws.Range("L132").Value = 7049.8329
$ws.Range("M132").Value = -351.8858
$ws.Range("N132").Value = -12109.8329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1979.579
$ws.Range("I134").Value = 1445.1395
$ws.Range("J134").Value = 3621.0715
$ws.Range("K134").Value = 4335.4185
$ws.Range("L134").Value = 10863.2145
$ws.Range("M134").Value = -1800.4185
$ws.Range("N134").Value = -15933.2145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1740.7797
$ws.Range("I31").Value = 1267.5135
$ws.Range("J31").Value = 2536.7273
$ws.Range("K31").Value = 1267.5135
$ws.Range("L31").Value = 2536.7273
$ws.Range("M31").Value = -972.5135
$ws.Range("N31").Value = -3126.7273
$ws.Range("H34").Value = 1740.7797
$ws.Range("I34").Value = 1267.5135
$ws.Range("J34").Value = 2536.7273
$ws.Range("K34").Value = 1267.5135
$ws.Range("L34").Value = 2536.7273
$ws.Range("M34").Value = -1065.5135
$ws.Range("N34").Value = -2940.7273
$ws.Range("H69").Value = 20095.555
$ws.Range("I69").Value = 10109.833
$ws.Range("J69").Value = 40067
$ws.Range("K69").Value = 10109.833
$ws.Range("L69").Value = 40067
$ws.Range("M69").Value = -9360.833000000001
$ws.Range("N69").Value = -41565
$ws.Range("H72").Value = 20095.555
$ws.Range("I72").Value = 10109.833
$ws.Range("J72").Value = 40067
$ws.Range("K72").Value = 30329.499
$ws.Range("L72").Value = 120201
$ws.Range("M72").Value = -26585.499
$ws.Range("N72").Value = -127689
$ws.Range("H132").Value = 356878.9
$ws.Range("I132").Value = 410443.7
$ws.Range("K132").Value = 1231331.1
$ws.Range("M132").Value = -1228801.1
$ws.Range("H134").Value = 1436.9807
$ws.Range("I134").Value = 1178.0555
$ws.Range("J134").Value = 2019.5625
$ws.Range("K134").Value = 3534.1665
$ws.Range("L134").Value = 6058.6875
$ws.Range("M134").Value = -999.1664999999998
$ws.Range("N134").Value = -11128.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3871.875
$ws.Range("I136").Value = 1099.8462
$ws.Range("J136").Value = 7147.909
$ws.Range("K136").Value = 3299.5386
$ws.Range("L136").Value = 21443.727
$ws.Range("M136").Value = 1800.4614
$ws.Range("N136").Value = -31643.727
$ws.Range("H139").Value = 2175.9714
$ws.Range("I139").Value = 1849.56
$ws.Range("J139").Value = 2992
$ws.Range("K139").Value = 5548.68
$ws.Range("L139").Value = 8976
$ws.Range("M139").Value = -408.6800000000003
$ws.Range("N139").Value = -19256
$ws.Range("H140").Value = 1591.2307
$ws.Range("I140").Value = 1333.56
$ws.Range("K140").Value = 4000.68
$ws.Range("M140").Value = 1179.32
$ws.Range("H141").Value = 3490.8667
$ws.Range("I141").Value = 3490.8667
$ws.Range("K141").Value = 10472.6001
$ws.Range("M141").Value = -5292.6001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6695.5
$ws.Range("I70").Value = 5800
$ws.Range("J70").Value = 7846.857
$ws.Range("K70").Value = 5800
$ws.Range("L70").Value = 7846.857
$ws.Range("M70").Value = -5530
$ws.Range("N70").Value = -8386.857
$ws.Range("H73").Value = 6695.5
$ws.Range("I73").Value = 5800
$ws.Range("J73").Value = 7846.857
$ws.Range("K73").Value = 5800
$ws.Range("L73").Value = 7846.857
$ws.Range("M73").Value = -4864
$ws.Range("N73").Value = -9718.857
$ws.Range("H122").Value = 2624.889
$ws.Range("I122").Value = 1972
$ws.Range("K122").Value = 5916
$ws.Range("M122").Value = -3466
$ws.Range("H132").Value = 1846.0238
$ws.Range("I132").Value = 1349.3846
$ws.Range("J132").Value = 2653.0625
$ws.Range("K132").Value = 4048.1538
$ws.Range("L132").Value = 7959.1875
$ws.Range("M132").Value = -1518.1538
$ws.Range("N132").Value = -13019.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 362.5
$ws.Range("I55").Value = 125
$ws.Range("J55").Value = 600
$ws.Range("K55").Value = 125
$ws.Range("L55").Value = 600
$ws.Range("M55").Value = 48
$ws.Range("N55").Value = -946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 800.38464
$ws.Range("I136").Value = 784.4138
$ws.Range("J136").Value = 846.7
$ws.Range("K136").Value = 2353.2414
$ws.Range("L136").Value = 2540.1
$ws.Range("M136").Value = 196.7586000000001
$ws.Range("N136").Value = -7640.1
